# Generate Report for Handback
#
# The localization-status workbook tracks handoff/handback state per
# target language (zh-cn, de-de). File "228b4934-...md" has just been
# handed back and is now in sync with en-US, so its Status moves from
# "Ready for handoff" to "Handed back: in sync with en-US" and its
# "Latest Handback DateTime" is stamped with the handback time, on both
# the zh-cn and de-de sheets (row 2 = the 228b4934 file in each table).
# The Overview sheet's per-language summary column for that same file
# (row 2, since A2 = 228b4934-...md there too) is updated to match.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_overview.Range("B2").Value = "Handed back: in sync with en-US"
$ws_overview.Range("C2").Value = "Handed back: in sync with en-US"

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("B2").Value = "Handed back: in sync with en-US"
$ws_zhcn.Range("G2").Value = "2016-02-22 14:03:21"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("B2").Value = "Handed back: in sync with en-US"
$ws_dede.Range("G2").Value = "2016-02-22 14:03:42"
